$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PAGOS")
$ws.Range("A1").Value = "test"
